# Update "Förändrad" (column C) date value from 2024-05-07 (45419) to
# 2024-05-08 (45420) for every existing data row (2-28).
$ws = $excel.ActiveWorkbook.ActiveSheet

for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45420
}

# Row 28 gains an explicit row height (matches the rest of the data rows).
$ws.Rows.Item(28).RowHeight = 15

# Append the new record as row 29.
$ws.Range("A29").Value2 = "A 17933-2024"

$ws.Range("B29").Value2 = 45419
$ws.Range("B29").NumberFormat = "YYYY-MM-DD"

$ws.Range("C29").Value2 = 45420
$ws.Range("C29").NumberFormat = "YYYY-MM-DD"

$ws.Range("D29").Value2 = "OKÄNT"
$ws.Range("E29").Value2 = "OKÄNT"

$ws.Range("G29").Value2 = 1.4
$ws.Range("H29").Value2 = 0
$ws.Range("I29").Value2 = 0
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 0
$ws.Range("L29").Value2 = 0
$ws.Range("M29").Value2 = 0
$ws.Range("N29").Value2 = 0
$ws.Range("O29").Value2 = 0
$ws.Range("P29").Value2 = 0
$ws.Range("Q29").Value2 = 0

$ws.Range("R29").WrapText = $true
$ws.Range("R29").Value2 = ""
